$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark (it currently sits in the
#    empty paragraph right after the Financial Updates paragraph).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Locate the run that holds the old sentence and replace it with
#    four separate runs (matching the target OOXML run-split), then
#    re-insert the "_GoBack" bookmark right after the new text, still
#    inside the same paragraph.
# ------------------------------------------------------------------
$find = $d.Content.Find
$find.Execute("Although this issue was resolved, it happens only during that timeframe, after rebooting the Raspberry Pi, the issues becomes unsolved again.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng = $find.Parent
$start = $rng.Start
$rng.Delete()

$r1 = $d.Range($start, $start)
$r1.InsertAfter("Although this issue was resolved, ")
$r1.LanguageID = "en-US"

$r2 = $d.Range($r1.End, $r1.End)
$r2.InsertAfter("I’ve noticed the RFID sensor was not soldered properly, this prevented the sensor from reading the tag properly, causing a big issue")
$r2.LanguageID = "en-US"

$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter(".")
$r3.LanguageID = "en-US"

$r4 = $d.Range($r3.End, $r3.End)
$r4.InsertAfter(" In the near future, this shall be fixed with better soldering.")
$r4.LanguageID = "en-US"

# Adding a bookmark with a collapsed range exactly at a paragraph-end
# boundary mis-anchors it, so insert a throwaway marker character,
# wrap the bookmark around it, then delete the marker -- the bookmark
# collapses back down to a zero-length span in the right spot.
$tmp = $d.Range($r4.End, $r4.End)
$tmp.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $tmp)
$tmp.Delete()
